$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: swap Exchange (Binance -> ByBit), fix From-date, swap Strategy (MACD -> Scalping1)
$ws.Range("B2").Value = "ByBit"
$ws.Range("D2").Value = 44197
$ws.Range("J2").Value = "Scalping1"

# Row 3: swap Exchange (ByBit -> Binance), fix From-date, swap Strategy (EarlyMACD -> Scalping1)
$ws.Range("B3").Value = "Binance"
$ws.Range("D3").Value = 44197
$ws.Range("J3").Value = "Scalping1"

# Row 4: remove the extra test case entirely (only the From/To date cells remain, blank)
$ws.Range("A4:C4").Clear()
$ws.Range("D4:E4").ClearContents()
$ws.Range("F4:J4").Clear()

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("J7").Select()
